$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new columns before the current column R (shifting old R:AE -> U:AH)
$ws.Range("R1:T1").EntireColumn.Insert()

# New header cells for the inserted columns
$ws.Range("R1").Value = "general_college_subjects.history"
$ws.Range("S1").Value = "general_college_subjects.electives"
$ws.Range("T1").Value = "general_college_subjects.cs"

# New data cells (row 2) for the inserted columns
$ws.Range("R2").Value = 2
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 0

# Normalize the "Considered"/"Not Considered" text to lowercase
$ws.Range("D2").Value = "considered"
$ws.Range("E2").Value = "considered"
$ws.Range("F2").Value = "not considered"
$ws.Range("G2").Value = "considered"
$ws.Range("H2").Value = "considered"
$ws.Range("I2").Value = "considered"
$ws.Range("J2").Value = "considered"
